$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.915.58'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '2.909.98'
$ws.Range('E3').Value = '  -4.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'585.91"
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').Value = "'146.34"
$ws.Range('E6').Value = '  -3.60%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.502"
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('D9').Value = '2.908.79'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').Value = "'6.73"
$ws.Range('E10').Value = '  +5.86%  '
$ws.Range('E11').Value = '  -5.15%  '
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('E13').Value = '  -4.53%  '
$ws.Range('D14').Value = "'33.53"
$ws.Range('E14').Value = '  -4.08%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '3.392.90'
$ws.Range('E16').Value = '  -4.19%  '
$ws.Range('D17').Value = '60.856.96'
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('D18').Value = "'6.77"
$ws.Range('E18').Value = '  -3.90%  '
$ws.Range('D19').Value = '2.909.53'
$ws.Range('E19').Value = '  -4.15%  '
$ws.Range('D20').Value = "'426.65"
$ws.Range('E20').Value = '  -6.57%  '
$ws.Range('E21').Value = '  -4.58%  '
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('E23').Value = '  -5.26%  '
$ws.Range('D24').Value = "'80.18"
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').Value = "'10.96"
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('E26').Value = '  -3.05%  '
$ws.Range('D27').Value = "'11.84"
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = "'7.23"
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('E31').Value = '  -3.52%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = "'26.44"
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('E34').Value = '  -4.68%  '
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('E37').Value = '  -4.91%  '
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('D39').Value = "'49.26"
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('D40').Value = "'2.01"
$ws.Range('E40').Value = '  -4.24%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  -4.96%  '
$ws.Range('E43').Value = '  -1.63%  '
$ws.Range('D44').Value = "'41.48"
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = "'376.50"
$ws.Range('E45').Value = '  -4.19%  '
$ws.Range('D46').Value = "'0.0346"
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('D47').Value = '2.669.37'
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('D48').Value = "'132.79"
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = "'24.71"
$ws.Range('E50').Value = '  +2.64%  '
$ws.Range('E51').Value = '  -1.74%  '
